$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New content -----------------------------------------------------
# Order matters: the shared-string table is built in first-write order,
# and the source workbook picked up "FillPageVehicleData" (row 15) before
# the new section xpaths on row 1, so mirror that sequence here.
$ws.Range("A15").Value2 = "103_TruckInsuranceAutomobile_001_SmokeTest_FillPageVehicleData"

$ws.Range("C1").Value2 = '//*[@id="insurance-form"]/div/section[1]'
$ws.Range("D1").Value2 = '//*[@id="insurance-form"]/div/section[2]'
$ws.Range("E1").Value2 = '//*[@id="insurance-form"]/div/section[3]'
$ws.Range("F1").Value2 = '//*[@id="insurance-form"]/div/section[4]'
$ws.Range("G1").Value2 = '//*[@id="insurance-form"]/div/section[5]'

# Row 2 stays empty in C:G but picks up a new yellow-fill style.
$ws.Range("C2:G2").Interior.Color = 65535

# --- Column widths -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 57.7213541666667
$ws.Range("C1:G1").EntireColumn.ColumnWidth = 34.8332992125984

# --- Picture reposition (keep the same pixel extent) --------------------
$shp = $ws.Shapes.Item(1)
$shp.Top = 331.8
$shp.Width = 1146.6066929133858
$shp.Height = 719.91

# --- Window / selection state -------------------------------------------
$ws.Range("K24").Select() | Out-Null
